# Append the new "Świeczki" product row (row 44) to the products sheet,
# extending the used range from A1:E43 to A1:E44.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(44, 1).Value = "P043"
$ws.Cells.Item(44, 2).Value = "Świeczki"
$ws.Cells.Item(44, 3).Value = "Dom"
$ws.Cells.Item(44, 4).Value = 20
$ws.Cells.Item(44, 5).Value = 100
